# Update the "Förändrad" (C) date column from 2023-09-15 (45184) to
# 2023-09-17 (45186) for every data row, and add the row's "Beteckning"
# (column A) as the friendly-name second argument of every HYPERLINK()
# formula in columns S, T, V, W, X, Y for rows 2-29 (the only rows that
# have those link formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used data row in column A (xlUp = -4162)
$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

# --- 1) Column C: 45184 -> 45186 for rows 2..lastRow ------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("C$r")
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45186
    }
}

# --- 2) Hyperlink formulas: add friendly name = column A value -------------
$linkCols = @{
    "S" = @("artfynd",          "xlsx")
    "T" = @("kartor",           "png")
    "V" = @("klagomål",         "docx")
    "W" = @("klagomålsmail",    "docx")
    "X" = @("tillsyn",          "docx")
    "Y" = @("tillsynsmail",     "docx")
}

for ($r = 2; $r -le $lastRow; $r++) {
    $beteckning = $ws.Range("A$r").Value2

    foreach ($col in $linkCols.Keys) {
        $linkCell = $ws.Range("$col$r")
        # Only touch cells that already hold a HYPERLINK(...) formula (single
        # URL argument, no friendly name yet) - this matches rows 2-29 in the
        # source workbook and leaves everything else untouched.
        if ($linkCell.HasFormula) {
            $folder = $linkCols[$col][0]
            $ext = $linkCols[$col][1]
            $url = "https://klasma.github.io/Logging_VALLENTUNA/$folder/$beteckning.$ext"
            $linkCell.Formula = "=HYPERLINK(""$url"", ""$beteckning"")"
        }
    }
}

Write-Host "Update complete"
